# myCBD/myInfo/raceLink.xlsx - "more standarization and SOPH prep"
#
# Sheet1:
#  - Column A (raceCode) values are re-standardized from single
#    letters (a,b,c,d,e,f,g,h,j) to short mnemonic codes
#    (AIAN,Black,Asian,Hisp,NHPI,White,Multi,Other,Unknown)
#  - Two new columns are appended: J "Ethan" and K "DCDC"
#  - A new Total row (row 12) is appended
#
# Sheet2 content is unchanged (only cosmetic/version metadata differs,
# which Excel manages automatically on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- New columns J (Ethan) and K (DCDC), header + body ----
$ws.Range("J1").Value = "Ethan"

$ws.Range("J4").Value = "ASIAN_NH"
$ws.Range("J3").Value = "BLACK_NH"
$ws.Range("J2").Value = "AIAN_NH"
$ws.Range("J7").Value = "WHITE_NH"
$ws.Range("J5").Value = "HISPANIC"
$ws.Range("J6").Value = "NHPI_NH"
$ws.Range("J8").Value = "MR_NH"

# ---- New Total row ----
$ws.Range("J12").Value = "TOTAL"
$ws.Range("A12").Value = "Total"
$ws.Range("F12").Value = "Total"

$ws.Range("K1").Value = "DCDC"
$ws.Range("K4").Value = "A"
$ws.Range("K7").Value = "W"
$ws.Range("K3").Value = "B"
$ws.Range("K5").Value = "H"
$ws.Range("K2").Value = "I"
$ws.Range("K8").Value = "M"
$ws.Range("K6").Value = "P"

# ---- Column A: standardize race codes ----
$ws.Range("A2").Value = "AIAN"
$ws.Range("A6").Value = "NHPI"
$ws.Range("A8").Value = "Multi"
$ws.Range("A3").Value = "Black"
$ws.Range("A4").Value = "Asian"
$ws.Range("A5").Value = "Hisp"
$ws.Range("A7").Value = "White"
$ws.Range("A9").Value = "Other"
$ws.Range("A10").Value = "Unknown"
$ws.Range("A11").Value = "Unknown"
